$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (shifts existing data rows 2-21 down to 3-22)
$ws.Rows.Item(2).Insert()
# Remove any formatting the insert carried over from the header row above
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row 2 with the new data values
$ws.Range("A2").Value = 0.1050096067542932
$ws.Range("B2").Value = -1.756468223065746
$ws.Range("C2").Value = 0.4945203567645989

# Append 9 brand-new rows of data after the (now) last row (row 22)
$newRows = @(
    @(-0.8633800915309378, -0.3154059344408438, 0.4824308418497782),
    @(-0.4081483519807154, -0.6726997543354425, -0.2190668820118418),
    @(0.2211332225373814, 0.241335413285664, 0.08368853798934378),
    @(0.06768137718341787, 0.3379019900244107, 0.1505034766635118),
    @(0.07254024853511698, 0.5556785336562575, -0.05807583201296457),
    @(0.1816357883567719, 0.1322741392923868, -0.08515337003128903),
    @(-0.02734556931013958, -0.1169588795425942, 0.04497027853313797),
    @(-0.02540700723017953, -0.06986615411481072, -0.074921377335808),
    @(0.02237761537639455, -0.07008743807863513, -0.003453258577050004)
)

$startRow = 23
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
}
